# Remove the trailing "Ver no Jupiter..." / copyright boilerplate block
# that followed the LOQ4095 requirement paragraph, along with the blank
# paragraph that separated them, while keeping the blank paragraph and
# page-break paragraph that come after the block.

$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter" text and the one
# that contains the copyright notice, then delete the range spanning from
# the blank paragraph right after the LOQ4095 line through the end of the
# copyright paragraph.

$jupiterPara = $null
$copyrightPara = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Ver no Jupiter*") {
        $jupiterPara = $p
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $copyrightPara = $p
    }
}

if ($jupiterPara -ne $null -and $copyrightPara -ne $null) {
    # The empty paragraph immediately preceding "Ver no Jupiter" is the one
    # that duplicated the blank separator left after the deletion, so start
    # the deletion there.
    $blankBefore = $jupiterPara.Previous()
    $startRange = $blankBefore.Range
    $endRange = $copyrightPara.Range

    $deleteRange = $d.Range($startRange.Start, $endRange.End)
    $deleteRange.Delete()
}
